$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1019
$ws.Range("F7").Value = 608
$ws.Range("F8").Value = 568
$ws.Range("F9").Value = 1482
$ws.Range("F10").Value = 138
$ws.Range("F11").Value = 1383
$ws.Range("F12").Value = 3032
$ws.Range("F13").Value = 504
$ws.Range("F14").Value = 1683
$ws.Range("F15").Value = 1371
$ws.Range("F16").Value = 818
$ws.Range("F17").Value = 251
$ws.Range("F18").Value = 1417
$ws.Range("F21").Value = 1153
$ws.Range("F22").Value = 21
$ws.Range("F23").Value = 413
$ws.Range("F24").Value = 30
$ws.Range("F25").Value = 3583
$ws.Range("F26").Value = 711
$ws.Range("F28").Value = 1581

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 39
$ws.Range("F7").Value = 5
$ws.Range("F13").Value = 87

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 20

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 20
$ws.Range("F8").Value = 39
$ws.Range("F12").Value = 5
$ws.Range("F16").Value = 1019
$ws.Range("F18").Value = 608
$ws.Range("F19").Value = 568
$ws.Range("F20").Value = 1482
$ws.Range("F21").Value = 138
$ws.Range("F22").Value = 1383
$ws.Range("F23").Value = 3032
$ws.Range("F24").Value = 504
$ws.Range("F25").Value = 1683
$ws.Range("F26").Value = 1371
$ws.Range("F27").Value = 818
$ws.Range("F28").Value = 251
$ws.Range("F29").Value = 1417
$ws.Range("F34").Value = 1153
$ws.Range("F35").Value = 21
$ws.Range("F36").Value = 413
$ws.Range("F37").Value = 30
$ws.Range("F38").Value = 3583
$ws.Range("F39").Value = 711
$ws.Range("F41").Value = 1581
$ws.Range("F42").Value = 87
